$wb = $excel.ActiveWorkbook

# --- Sheet "Resumen": update zone label and max value ---
$wsResumen = $wb.Worksheets.Item("Resumen")
$wsResumen.Range("B2").Value = "Z2"
$wsResumen.Range("C2").Value = 556.2911287706779

# --- Sheet "Solucion": final (shuffled) Pedido/Salida assignment ---
$wsSolucion = $wb.Worksheets.Item("Solucion")

$solucionData = @(
    @("Pedido_44", "S001"),
    @("Pedido_50", "S021"),
    @("Pedido_7", "S041"),
    @("Pedido_60", "S031"),
    @("Pedido_37", "S011"),
    @("Pedido_52", "S051"),
    @("Pedido_20", "S022"),
    @("Pedido_14", "S002"),
    @("Pedido_39", "S042"),
    @("Pedido_4", "S032"),
    @("Pedido_16", "S012"),
    @("Pedido_28", "S052"),
    @("Pedido_34", "S023"),
    @("Pedido_58", "S043"),
    @("Pedido_27", "S003"),
    @("Pedido_41", "S033"),
    @("Pedido_53", "S053"),
    @("Pedido_8", "S044"),
    @("Pedido_25", "S013"),
    @("Pedido_5", "S024"),
    @("Pedido_32", "S054"),
    @("Pedido_45", "S004"),
    @("Pedido_2", "S034"),
    @("Pedido_31", "S045"),
    @("Pedido_59", "S014"),
    @("Pedido_10", "S025"),
    @("Pedido_18", "S055"),
    @("Pedido_38", "S005"),
    @("Pedido_6", "S046"),
    @("Pedido_29", "S035"),
    @("Pedido_42", "S015"),
    @("Pedido_48", "S006"),
    @("Pedido_33", "S026"),
    @("Pedido_35", "S056"),
    @("Pedido_19", "S036"),
    @("Pedido_56", "S047"),
    @("Pedido_43", "S016"),
    @("Pedido_13", "S027"),
    @("Pedido_54", "S007"),
    @("Pedido_3", "S057"),
    @("Pedido_9", "S037"),
    @("Pedido_24", "S048"),
    @("Pedido_40", "S017"),
    @("Pedido_57", "S028"),
    @("Pedido_49", "S008"),
    @("Pedido_21", "S058"),
    @("Pedido_17", "S038"),
    @("Pedido_1", "S049"),
    @("Pedido_51", "S018"),
    @("Pedido_23", "S029"),
    @("Pedido_55", "S059"),
    @("Pedido_11", "S009"),
    @("Pedido_26", "S050"),
    @("Pedido_15", "S039"),
    @("Pedido_47", "S019"),
    @("Pedido_12", "S030"),
    @("Pedido_46", "S010"),
    @("Pedido_22", "S060"),
    @("Pedido_30", "S040"),
    @("Pedido_36", "S020")
)

for ($i = 0; $i -lt $solucionData.Length; $i++) {
    $row = $i + 2
    $wsSolucion.Cells.Item($row, 1).Value = $solucionData[$i][0]
    $wsSolucion.Cells.Item($row, 2).Value = $solucionData[$i][1]
}

# --- Sheet "Metricas": update per-zone time values ---
$wsMetricas = $wb.Worksheets.Item("Metricas")
$wsMetricas.Range("B2").Value = 547.0383392799222
$wsMetricas.Range("B3").Value = 556.2911287706779
$wsMetricas.Range("B4").Value = 554.3032544058817
